$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "June_Planner" (first sheet) - fill in the Activity-1/2/3 Status
# columns (D/E/F) for rows 9-12 with "Complete" / "In Progress" markers,
# reusing the already-present formatting (style index 21 = "Complete" green
# font, style index 22 = "In Progress" accent font) from row 2 as the
# formatting source so the new cells pick up identical styling.
# ---------------------------------------------------------------------------
$planner = $wb.Worksheets.Item("June_Planner")
$planner.Activate()

$completeFmt = $planner.Range("D2")   # style s="21" -> "Complete"
$inProgressFmt = $planner.Range("F2") # style s="22" -> "In Progress"

# Row 9: Complete, Complete, Complete
$completeFmt.Copy()
$planner.Range("D9").PasteSpecial(-4122)
$planner.Range("D9").Value = "Complete"
$completeFmt.Copy()
$planner.Range("E9").PasteSpecial(-4122)
$planner.Range("E9").Value = "Complete"
$completeFmt.Copy()
$planner.Range("F9").PasteSpecial(-4122)
$planner.Range("F9").Value = "Complete"

# Row 10: Complete, In Progress, Complete
$completeFmt.Copy()
$planner.Range("D10").PasteSpecial(-4122)
$planner.Range("D10").Value = "Complete"
$inProgressFmt.Copy()
$planner.Range("E10").PasteSpecial(-4122)
$planner.Range("E10").Value = "In Progress"
$completeFmt.Copy()
$planner.Range("F10").PasteSpecial(-4122)
$planner.Range("F10").Value = "Complete"

# Row 11: Complete, Complete, In Progress
$completeFmt.Copy()
$planner.Range("D11").PasteSpecial(-4122)
$planner.Range("D11").Value = "Complete"
$completeFmt.Copy()
$planner.Range("E11").PasteSpecial(-4122)
$planner.Range("E11").Value = "Complete"
$inProgressFmt.Copy()
$planner.Range("F11").PasteSpecial(-4122)
$planner.Range("F11").Value = "In Progress"

# Row 12: In Progress, Complete, Complete
$inProgressFmt.Copy()
$planner.Range("D12").PasteSpecial(-4122)
$planner.Range("D12").Value = "In Progress"
$completeFmt.Copy()
$planner.Range("E12").PasteSpecial(-4122)
$planner.Range("E12").Value = "Complete"
$completeFmt.Copy()
$planner.Range("F12").PasteSpecial(-4122)
$planner.Range("F12").Value = "Complete"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet "Sheet1" (third sheet) - same Complete/In Progress pattern for the
# previously-blank D/E/F columns of rows 2-6 (row 6 only gets D6 filled in,
# matching the target state).
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()

# Row 2: Complete, Complete, Complete
$completeFmt.Copy()
$sheet1.Range("D2").PasteSpecial(-4122)
$sheet1.Range("D2").Value = "Complete"
$completeFmt.Copy()
$sheet1.Range("E2").PasteSpecial(-4122)
$sheet1.Range("E2").Value = "Complete"
$completeFmt.Copy()
$sheet1.Range("F2").PasteSpecial(-4122)
$sheet1.Range("F2").Value = "Complete"

# Row 3: Complete, In Progress, Complete
$completeFmt.Copy()
$sheet1.Range("D3").PasteSpecial(-4122)
$sheet1.Range("D3").Value = "Complete"
$inProgressFmt.Copy()
$sheet1.Range("E3").PasteSpecial(-4122)
$sheet1.Range("E3").Value = "In Progress"
$completeFmt.Copy()
$sheet1.Range("F3").PasteSpecial(-4122)
$sheet1.Range("F3").Value = "Complete"

# Row 4: Complete, Complete, In Progress
$completeFmt.Copy()
$sheet1.Range("D4").PasteSpecial(-4122)
$sheet1.Range("D4").Value = "Complete"
$completeFmt.Copy()
$sheet1.Range("E4").PasteSpecial(-4122)
$sheet1.Range("E4").Value = "Complete"
$inProgressFmt.Copy()
$sheet1.Range("F4").PasteSpecial(-4122)
$sheet1.Range("F4").Value = "In Progress"

# Row 5: In Progress, Complete, Complete
$inProgressFmt.Copy()
$sheet1.Range("D5").PasteSpecial(-4122)
$sheet1.Range("D5").Value = "In Progress"
$completeFmt.Copy()
$sheet1.Range("E5").PasteSpecial(-4122)
$sheet1.Range("E5").Value = "Complete"
$completeFmt.Copy()
$sheet1.Range("F5").PasteSpecial(-4122)
$sheet1.Range("F5").Value = "Complete"

# Row 6: In Progress (D only; E6/F6 stay blank)
$inProgressFmt.Copy()
$sheet1.Range("D6").PasteSpecial(-4122)
$sheet1.Range("D6").Value = "In Progress"

$excel.CutCopyMode = $false

# Record the UI selection left on Sheet1 ...
$sheet1.Range("M3").Select()

# ... then re-activate June_Planner with C12 selected so it ends up as the
# workbook's active/visible tab, matching the target state.
$planner.Activate()
$planner.Range("C12").Select()
